$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update existing rows 302-307 (revised Imacec series) ---
$ws.Range("B302").Value = 106.7
$ws.Range("C302").Value = 107.1
$ws.Range("D302").Value = 95.8
$ws.Range("E302").Value = 104.3
$ws.Range("F302").Value = 120.3
$ws.Range("G302").Value = 129.5

$ws.Range("B303").Value = 104.2
$ws.Range("C303").Value = 104.5
$ws.Range("D303").Value = 88.1
$ws.Range("E303").Value = 101.7
$ws.Range("F303").Value = 122.8
$ws.Range("G303").Value = 126.4
$ws.Range("H303").Value = 98.1
$ws.Range("I303").Value = 103.5
$ws.Range("J303").Value = 106.3

$ws.Range("B304").Value = 120.5
$ws.Range("C304").Value = 115.1
$ws.Range("D304").Value = 101.4
$ws.Range("E304").Value = 114
$ws.Range("F304").Value = 129.2
$ws.Range("G304").Value = 143
$ws.Range("H304").Value = 118.9
$ws.Range("I304").Value = 119.8
$ws.Range("J304").Value = 123

$ws.Range("B305").Value = 112.8
$ws.Range("C305").Value = 103.2
$ws.Range("D305").Value = 98.3
$ws.Range("E305").Value = 107.5
$ws.Range("F305").Value = 104.6
$ws.Range("G305").Value = 124.3
$ws.Range("I305").Value = 112.8
$ws.Range("J305").Value = 114.7

$ws.Range("B306").Value = 116.6
$ws.Range("C306").Value = 103.7
$ws.Range("D306").Value = 100.7
$ws.Range("E306").Value = 108.6
$ws.Range("F306").Value = 102.5
$ws.Range("G306").Value = 134.6
$ws.Range("H306").Value = 120.6
$ws.Range("I306").Value = 115.3
$ws.Range("J306").Value = 118.6

$ws.Range("B307").Value = 115.9
$ws.Range("C307").Value = 101.2
$ws.Range("D307").Value = 100.2
$ws.Range("E307").Value = 108.4
$ws.Range("F307").Value = 96.3
$ws.Range("G307").Value = 130.9
$ws.Range("H307").Value = 122.2
$ws.Range("I307").Value = 114.8
$ws.Range("J307").Value = 118

# --- Append new row 308 (01-07-2021) ---
$ws.Range("A308").NumberFormat = "@"
$ws.Range("A308").Value = "01-07-2021"
$ws.Range("A308").Style = "Normal"
$ws.Range("B308").Value = 113.5
$ws.Range("C308").Value = 100.6
$ws.Range("D308").Value = 99.9
$ws.Range("E308").Value = 107.5
$ws.Range("F308").Value = 95.7
$ws.Range("G308").Value = 135.3
$ws.Range("H308").Value = 115.3
$ws.Range("I308").Value = 111.5
$ws.Range("J308").Value = 115.2
